# Insert a new record row at row 551 (pushing existing rows 551:614 down to
# 552:615) and populate it with the new observation's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("551:551").Insert()

$ws.Range("A551").Value = 3
$ws.Range("B551").Value = 'Femacal de La Calera'
$ws.Range("C551").Value = 'Coquimbo'
$ws.Range("D551").Value = 44946
$ws.Range("E551").Value = 5
$ws.Range("F551").Value = 100112021
$ws.Range("G551").Value = 'Ají'
$ws.Range("H551").Value = 'Inferno'
$ws.Range("I551").Value = 'Primera'
$ws.Range("J551").Value = 73
$ws.Range("K551").Value = 18000
$ws.Range("L551").Value = 19000
$ws.Range("M551").Value = 18479
$ws.Range("N551").Value = '$/caja 15 kilos'
$ws.Range("O551").Value = 'Limache'
$ws.Range("P551").Value = 1232
$ws.Range("Q551").Value = 15
$ws.Range("R551").Value = 'Hortaliza'
